# "Copy and Delete sheet activities"
# 1. Activate Sheet1 and select H9 on it.
# 2. Copy Sheet1 to the end of the workbook and rename the copy to "CopiedSheet".
# 3. Clear the copied sheet's sort state (the copy inherited Sheet1's sortState
#    for A2:G12, which no longer matches once a row is removed below).
# 4. Delete the "Boris" row (row 10) from the copied sheet.
# The copy activity carries over the active selection, so both Sheet1 and the
# new CopiedSheet end up with the H9 selection; CopiedSheet becomes the active
# (tabSelected) sheet, and Sheet3 (previously active) loses tabSelected.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$sheet1.Range("H9").Select() | Out-Null

$sheet1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "CopiedSheet"

$newSheet.Sort.SortFields.Clear()
$newSheet.Rows.Item(10).Delete()
